# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1. Metadata sheet: bump the generation Date.
# 2. Elements sheet: the two "Mapping" columns (AK = RIM Mapping, AL = the
#    new business mapping) were swapped so the new
#    "Spécification métier" mapping comes first (column AK) and the RIM
#    mapping moves to column AL - header text, column width and every data
#    row all move together.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 -> new Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) -----------------------
$ws = $wb.Worksheets.Item("Elements")

# Header row (row 1)
$ws.Range("AK1").Value = "Mapping: Spécification métier vers l'extension ROR RORPractitionerRoleName"
$ws.Range("AL1").Value = "Mapping: RIM Mapping"

# Column widths follow their column, so they swap too (best-effort - COM
# ColumnWidth works in whole/sixth character units so this lands on the
# closest value Excel itself would store).
$ws.Columns.Item(37).ColumnWidth = 82.18359375
$ws.Columns.Item(38).ColumnWidth = 24.98046875

# Data rows 2-22: new AK value = old AL value, new AL value = old AK value.
$rows = @(
  @{Row=2;  AK='';                  AL=''},
  @{Row=3;  AK='';                  AL='n/a'},
  @{Row=4;  AK='';                  AL=''},
  @{Row=5;  AK='';                  AL=''},
  @{Row=6;  AK='';                  AL='n/a'},
  @{Row=7;  AK='';                  AL=''},
  @{Row=8;  AK='';                  AL='N/A'},
  @{Row=9;  AK='';                  AL='N/A'},
  @{Row=10; AK='civiliteExercice';  AL='N/A'},
  @{Row=11; AK='';                  AL=''},
  @{Row=12; AK='';                  AL='n/a'},
  @{Row=13; AK='';                  AL=''},
  @{Row=14; AK='';                  AL='N/A'},
  @{Row=15; AK='nomExercice';       AL='N/A'},
  @{Row=16; AK='';                  AL=''},
  @{Row=17; AK='';                  AL='n/a'},
  @{Row=18; AK='';                  AL=''},
  @{Row=19; AK='';                  AL='N/A'},
  @{Row=20; AK='prenomExercice';    AL='N/A'},
  @{Row=21; AK='';                  AL='N/A'},
  @{Row=22; AK='';                  AL='N/A'}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 37).Value = $r.AK
    $ws.Cells.Item($r.Row, 38).Value = $r.AL
}
